$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.589.91"
$ws.Range("E2").Value = "  +0.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.797.82"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.77%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "340.01"
$ws.Range("E5").Value = "  +1.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9977"
$ws.Range("E6").Value = "  -0.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3927"
$ws.Range("E7").Value = "  +3.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3467"
$ws.Range("E8").Value = "  -0.50%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.20"
$ws.Range("E9").Value = "  -1.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.197"
$ws.Range("E10").Value = "  -1.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07513"
$ws.Range("E11").Value = "  -0.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9960"
$ws.Range("E12").Value = "  -0.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.98"
$ws.Range("E13").Value = "  -0.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.521"
$ws.Range("E14").Value = "  -0.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.796.13"
$ws.Range("E15").Value = "  -0.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.154"
$ws.Range("E16").Value = "  +0.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001102"
$ws.Range("E17").Value = "  -0.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06722"
$ws.Range("E18").Value = "  +0.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.07"
$ws.Range("E19").Value = "  -0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9971"
$ws.Range("E20").Value = "  -0.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.74"
$ws.Range("E21").Value = "  +1.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.563"
$ws.Range("E22").Value = "  +1.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.567.86"
$ws.Range("E23").Value = "  +0.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.45"
$ws.Range("E24").Value = "  -1.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.403"
$ws.Range("E25").Value = "  -2.27%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.29"
$ws.Range("E26").Value = "  -1.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.511"
$ws.Range("E27").Value = "  -2.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.468"
$ws.Range("E28").Value = "  +0.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "155.81"
$ws.Range("E29").Value = "  +3.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.998.03"
$ws.Range("E30").Value = "  -0.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "135.36"
$ws.Range("E31").Value = "  +0.94%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.038"
$ws.Range("E32").Value = "  -1.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.051"
$ws.Range("E33").Value = "  -1.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08782"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.12"
$ws.Range("E35").Value = "  -1.77%  "

$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.469"
$ws.Range("E36").Value = "  -0.92%  "

$ws.Range("B37").Value = "WEMIXTOKEN"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.615"
$ws.Range("E37").Value = "  -3.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02418"
$ws.Range("E38").Value = "  +2.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06473"
$ws.Range("E39").Value = "  +1.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6837"
$ws.Range("E40").Value = "  -0.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2215"
$ws.Range("E41").Value = "  -0.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.252"
$ws.Range("E42").Value = "  -2.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.380"
$ws.Range("E43").Value = "  -5.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.64"
$ws.Range("E44").Value = "  +2.07%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6428"
$ws.Range("E45").Value = "  -0.49%  "

$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9964"
$ws.Range("E46").Value = "  -0.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.873"
$ws.Range("E47").Value = "  +0.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.144"
$ws.Range("E48").Value = "  +0.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.31"
$ws.Range("E49").Value = "  +0.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07207"
$ws.Range("E50").Value = "  -0.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.98"
$ws.Range("E51").Value = "  +0.19%  "
